$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post at row 515 ("「一緒に読むのが好き」") was removed from the source data.
# Deleting the entire row shifts every subsequent row up by one and
# automatically shrinks the sheet's used range (dimension) from C684 to C683.
$ws.Rows("515:515").Delete()
